# Refresh the cryptos price list (Coin/Link/Price/Volume columns) to the
# latest scrape, matching the "Updated cryptos list ... with GitHub Actions"
# commit. Column A (rank index) and the header row are untouched; only
# B:E for data rows 2-51 are updated where the scrape produced new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range('D2').Value = '37.458.38'
$ws.Range('E2').Value = '  +1.22%  '

$ws.Range('D3').Value = '2.031.74'
$ws.Range('E3').Value = '  +0.65%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '''229.15'
$ws.Range('E5').Value = '  +1.47%  '

$ws.Range('D6').Value = '''0.613'
$ws.Range('E6').Value = '  +0.89%  '

$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').Value = '''55.89'
$ws.Range('E8').Value = '  +1.89%  '

$ws.Range('D9').Value = '''0.381'
$ws.Range('E9').Value = '  +0.30%  '

$ws.Range('D10').Value = '''0.0799'
$ws.Range('E10').Value = '  +1.61%  '

$ws.Range('E11').Value = '  -1.01%  '

$ws.Range('D12').Value = '2.329.54'
$ws.Range('E12').Value = '  +0.51%  '

$ws.Range('D13').Value = '''14.35'
$ws.Range('E13').Value = '  +0.58%  '

$ws.Range('D14').Value = '''20.19'
$ws.Range('E14').Value = '  -0.73%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''5.21'
$ws.Range('E15').Value = '  +1.49%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.741'
$ws.Range('E16').Value = '  -0.07%  '

$ws.Range('D17').Value = '2.028.83'
$ws.Range('E17').Value = '  +0.72%  '

$ws.Range('D18').Value = '37.367.44'
$ws.Range('E18').Value = '  +1.32%  '

$ws.Range('D19').Value = '''6.18'
$ws.Range('E19').Value = '  -0.83%  '

$ws.Range('D20').Value = '''69.00'

$ws.Range('D21').Value = '0.0₃0824'
$ws.Range('E21').Value = '  +0.43%  '

$ws.Range('D22').Value = '''223.05'
$ws.Range('E22').Value = '  -1.25%  '

$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  -0.01%  '

$ws.Range('E24').Value = '  +1.70%  '

$ws.Range('E25').Value = '  +2.97%  '

$ws.Range('D26').Value = '''164.96'
$ws.Range('E26').Value = '  -0.38%  '

$ws.Range('E27').Value = '  -1.28%  '

$ws.Range('E28').Value = '  +4.43%  '

$ws.Range('D29').Value = '''18.76'
$ws.Range('E29').Value = '  +0.40%  '

$ws.Range('E30').Value = '  -1.17%  '

$ws.Range('E31').Value = '  +0.56%  '

$ws.Range('D32').Value = '''4.48'
$ws.Range('E32').Value = '  -0.01%  '

$ws.Range('D33').Value = '''0.0605'
$ws.Range('E33').Value = '  -1.89%  '

$ws.Range('E34').Value = '  +0.84%  '

$ws.Range('E35').Value = '  +8.37%  '

$ws.Range('D36').Value = '''2.32'
$ws.Range('E36').Value = '  -1.04%  '

$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').Value = '''5.75'
$ws.Range('E37').Value = '  +8.66%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''3.24'
$ws.Range('E38').Value = '  +2.66%  '

$ws.Range('E39').Value = '  +0.09%  '

$ws.Range('D40').Value = '1.474.76'
$ws.Range('E40').Value = '  -0.64%  '

$ws.Range('E41').Value = '  -1.69%  '

$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '''2.84'
$ws.Range('E42').Value = '  +3.31%  '

$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').Value = '''0.0932'
$ws.Range('E43').Value = '  +0.51%  '

$ws.Range('D44').Value = '''95.05'
$ws.Range('E44').Value = '  -0.33%  '

$ws.Range('D45').Value = '''4.27'
$ws.Range('E45').Value = '  +18.58%  '

$ws.Range('D46').Value = '''16.31'
$ws.Range('E46').Value = '  -5.16%  '

$ws.Range('E47').Value = '  -2.63%  '

$ws.Range('D48').Value = '''1.01'
$ws.Range('E48').Value = '  +0.25%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '''7.10'
$ws.Range('E49').Value = '  -3.49%  '

$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = '''2.94'
$ws.Range('E50').Value = '  +0.56%  '

$ws.Range('D51').Value = '2.219.78'
$ws.Range('E51').Value = '  +0.49%  '
